$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2022-08-28 07:01:30"
$newTimestamp = "2022-08-28 20:57:18"

for ($row = 2; $row -le 35; $row++) {
    $cell = $ws.Cells.Item($row, 15)  # Column O is the 15th column
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
